$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'26.563.73"
$ws.Range("E2").Value = "  +0.47%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.848.26"
$ws.Range("E3").Value = "  +0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'264.44"
$ws.Range("E5").Value = "  +1.47%  "

# Row 6 - USDC
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.13%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.5249"
$ws.Range("E7").Value = "  +0.61%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3244"
$ws.Range("E8").Value = "  +0.36%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.06813"
$ws.Range("E9").Value = "  +0.81%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'18.92"
$ws.Range("E10").Value = "  +0.32%  "

# Row 11 - Polygon
$ws.Range("D11").Value = "'0.7818"
$ws.Range("E11").Value = "  +1.53%  "

# Row 12 - TRON
$ws.Range("D12").Value = "'0.07788"
$ws.Range("E12").Value = "  +0.72%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "'1.856.88"
$ws.Range("E13").Value = "  -0.11%  "

# Row 14 - Litecoin
$ws.Range("E14").Value = "  -0.27%  "

# Row 15 - Polkadot
$ws.Range("D15").Value = "'5.026"
$ws.Range("E15").Value = "  -0.13%  "

# Row 16 - BinanceUSD
$ws.Range("E16").Value = "  -0.08%  "

# Row 17 - Avalanche
$ws.Range("D17").Value = "'14.00"
$ws.Range("E17").Value = "  -1.03%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "'0.000007993"
$ws.Range("E18").Value = "  +1.04%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.02%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "'26.576.35"
$ws.Range("E20").Value = "  +0.26%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'4.641"
$ws.Range("E21").Value = "  +2.30%  "

# Row 22 - Cosmos
$ws.Range("D22").Value = "'9.491"
$ws.Range("E22").Value = "  -0.38%  "

# Row 23 - Chainlink
$ws.Range("D23").Value = "'6.027"
$ws.Range("E23").Value = "  +1.70%  "

# Row 24 - Monero
$ws.Range("D24").Value = "'142.91"
$ws.Range("E24").Value = "  -0.76%  "

# Row 25 - LidoDAOToken
$ws.Range("D25").Value = "'2.192"
$ws.Range("E25").Value = "  -6.92%  "

# Row 26 - Toncoin
$ws.Range("D26").Value = "'1.685"
$ws.Range("E26").Value = "  +1.56%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'17.04"
$ws.Range("E27").Value = "  +0.34%  "

# Row 28 - BitcoinCash
$ws.Range("D28").Value = "'111.72"
$ws.Range("E28").Value = "  +0.13%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "'4.198"
$ws.Range("E29").Value = "  +0.21%  "

# Row 30 - Filecoin
$ws.Range("D30").Value = "'4.125"
$ws.Range("E30").Value = "  -0.97%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "'0.08734"
$ws.Range("E31").Value = "  -0.57%  "

# Row 32 - Hedera
$ws.Range("E32").Value = "  +0.93%  "

# Row 33 - ARBITRUM
$ws.Range("D33").Value = "'1.135"
$ws.Range("E33").Value = "  +0.08%  "

# Row 34 - ImmutableX
$ws.Range("D34").Value = "'0.7230"
$ws.Range("E34").Value = "  +4.56%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "'2.879"
$ws.Range("E35").Value = "  +0.68%  "

# Row 36 - MXToken
$ws.Range("D36").Value = "'3.107"
$ws.Range("E36").Value = "  -0.18%  "

# Row 37 - RenderToken
$ws.Range("D37").Value = "'2.270"
$ws.Range("E37").Value = "  +2.76%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "'0.01793"
$ws.Range("E38").Value = "  +0.03%  "

# Row 39 - TheSandbox
$ws.Range("D39").Value = "'0.4874"
$ws.Range("E39").Value = "  -0.79%  "

# Row 40 - TrustWalletToken
$ws.Range("D40").Value = "'0.9027"
$ws.Range("E40").Value = "  +0.60%  "

# Row 41 - Quant
$ws.Range("D41").Value = "'110.89"
$ws.Range("E41").Value = "  -1.95%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "'5.991"
$ws.Range("E42").Value = "  -3.32%  "

# Row 43 - PaxDollar
$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  -0.14%  "

# Row 44 - Aptos
$ws.Range("D44").Value = "'7.667"
$ws.Range("E44").Value = "  -1.68%  "

# Row 45 - Decentraland
$ws.Range("D45").Value = "'0.4215"
$ws.Range("E45").Value = "  +0.27%  "

# Row 46 - Cronos
$ws.Range("D46").Value = "'0.05889"
$ws.Range("E46").Value = "  +0.10%  "

# Row 47 - EnergySwap
$ws.Range("D47").Value = "'9.065"
$ws.Range("E47").Value = "  -0.37%  "

# Row 48 - was Elrond, now Algorand
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1237"
$ws.Range("E48").Value = "  -2.11%  "

# Row 49 - was Algorand, now Elrond
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'35.07"
$ws.Range("E49").Value = "  -0.95%  "

# Row 50 - EOS
$ws.Range("D50").Value = "'0.8898"
$ws.Range("E50").Value = "  +3.68%  "

# Row 51 - Aave
$ws.Range("D51").Value = "'60.07"
$ws.Range("E51").Value = "  +1.47%  "
